$wb = $excel.ActiveWorkbook
$sheetIndexes = @(1, 4)
foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("F2").Value = 3143
    $ws.Range("F3").Value = 192
    $ws.Range("F6").Value = 1746
    $ws.Range("F8").Value = 99
    $ws.Range("F10").Value = 8
    $ws.Range("F11").Value = 1443
    $ws.Range("F13").Value = 571
    $ws.Range("F14").Value = 358
    $ws.Range("F15").Value = 75
    $ws.Range("F17").Value = 81
    $ws.Range("F21").Value = 97
    $ws.Range("F23").Value = 3395
    $ws.Range("F24").Value = 410
    $ws.Range("F25").Value = 294
    $ws.Range("F26").Value = 467
    $ws.Range("F27").Value = 59
    $ws.Range("F30").Value = 1132
    $ws.Range("F31").Value = 119
}
